# The "school_id" column (G) held a handful of generic "RRSxxxxxXX" codes
# that were reused across several different schools/teachers. Split those
# out into one unique Reading-Recovery school code per teacher/school row
# (commit: "couldnt find a way to nest the 3 dfs in pandas dataclass, so
# created one with school, teacher as index").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @("G2",  "RRS1080080"),
    @("G3",  "RRS1080080"),
    @("G4",  "RRS2030220"),
    @("G5",  "RRS2030220"),
    @("G6",  "RRS2030220"),
    @("G7",  "RRS2030220"),
    @("G8",  "RRS2010080"),
    @("G9",  "RRS4030010"),
    @("G10", "RRS4030010"),
    @("G11", "RRS4030010"),
    @("G12", "RRS4030010"),
    @("G13", "RRS2020030"),
    @("G14", "RRS2020030"),
    @("G15", "RRS2020030"),
    @("G16", "RRS2020030"),
    @("G17", "RRS2070510"),
    @("G18", "RRS2070510"),
    @("G19", "RRS2070510"),
    @("G20", "RRS2070510"),
    @("G21", "RRS2020080"),
    @("G22", "RRS2020080"),
    @("G23", "RRS2020080"),
    @("G24", "RRS2010450"),
    @("G25", "RRS2010450"),
    @("G26", "RRS2070140"),
    @("G27", "RRS2070140"),
    @("G28", "RRS2070140"),
    @("G29", "RRS2070140"),
    @("G30", "RRS2030250"),
    @("G31", "RRS2030250"),
    @("G32", "RRS2030250"),
    @("G33", "RRS2030250"),
    @("G34", "RRS2030250")
)

foreach ($pair in $updates) {
    $ws.Range($pair[0]).Value = $pair[1]
}

# Matches the final cursor position recorded in the saved workbook.
$ws.Range("G2").Select()
